$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 9 (shifts old rows 9-16 down to 10-17)
$ws.Rows.Item(9).Insert()

# Fill in the new row 7 with Address / adr
$ws.Range("B7").Value = "Address"
$ws.Range("C7").Value = "adr"

# Set the active selection to C8 to match sheetView selection
$ws.Range("C8").Select()
